$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header date labels (shift one month forward) ---
# Assigning a plain string like "June 2024"/"July 2024" directly to .Value
# gets auto-parsed by Excel as a date and reformatted with a date number
# format, which would introduce an unwanted style change. Instead, enter
# each value as a text formula (quoted literal) and then convert it to a
# static value via Copy / Paste-Special-Values, exactly like using Excel's
# "Paste Values" button. This keeps the cells as plain shared-string text
# cells with no style applied, matching the original formatting.
$ws.Range("A1").Formula = "=""June 2024"""
$ws.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4163)

$ws.Range("G1").Formula = "=""July 2024"""
$ws.Range("G1").Copy()
$ws.Range("G1").PasteSpecial(-4163)

# --- Update the waterfall data values in row 2 ---
$ws.Range("A2").Value = 1.502
$ws.Range("B2").Value = -0.013
$ws.Range("C2").Value = -0.028
$ws.Range("D2").Value = 0.079
$ws.Range("E2").Value = -0.03
$ws.Range("F2").Value = 0.043
$ws.Range("G2").Value = 1.549
